# Regenerate the localization handoff report: the rows that were just
# re-handed-off get their Priority set to "ht" (handoff type) and their
# timestamps bumped to reflect the new run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$rows = @(7, 8, 9, 10, 12, 13)

foreach ($r in $rows) {
    # Priority column (E) becomes "ht" on both language sheets.
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"

    # Latest HO Xliff Generate Date on the Overview sheet.
    $overview.Range("G$r").Value = "2016-08-24 12:22:19"

    # Latest Handoff Datetime on the zh-cn sheet.
    $zhcn.Range("H$r").Value = "2016-08-24 12:22:13"
}
